# Progress update as of 04-Nov-2025:
#  - Column H ("PERIOD TO EXPIRE") drops by 1 day for every training row.
#  - Column I ("LAST UPDATE") moves from 03-Nov-2025 to 04-Nov-2025.
#
# Column I holds a literal text date string (not a real Excel date), so we
# can't just assign the string to .Value - Excel's smart-entry would parse
# "04-Nov-2025" as an actual date serial and reformat the cell. Instead we
# stage the literal text via a formula in a scratch cell (which forces a
# text result), copy it, and paste-special "values only" into each target
# cell - this preserves the original General number format/style and keeps
# the stored type as plain text, just like the source cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$scratch = $ws.Range("Z1")
$xlPasteValues = -4163

for ($row = 3; $row -le 24; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H: PERIOD TO EXPIRE
    $hCell.Value = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I: LAST UPDATE
    $scratch.Formula = "=""04-Nov-2025"""
    $scratch.Copy()
    $iCell.PasteSpecial($xlPasteValues)
}

$scratch.Clear()
$excel.CutCopyMode = 0
